# Auto-generated Excel COM-interop script to apply daily data update
# for paises.xlsx ("Update countries & provincias Spain")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last refreshed" timestamp in the title cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 28 de Marzo de 2020 a las 20:29"

# Update per-row country data (column A holds the country name; the sheet is
# kept sorted descending by column B "Casos totales", so as totals change the
# country occupying a given row, and its stats in B:H, change too).

# Row 32
$ws.Range("A32").Value = "Polonia"
$ws.Range("B32").Value = 1638
$ws.Range("C32").Value = 249
$ws.Range("D32").Value = 7
$ws.Range("E32").Value = 1613
$ws.Range("F32").Value = 3
$ws.Range("G32").Value = 2
$ws.Range("H32").Value = 18

# Row 33
$ws.Range("A33").Value = "Japon"
$ws.Range("B33").Value = 1499
$ws.Range("C33").Value = 0
$ws.Range("D33").Value = 404
$ws.Range("E33").Value = 1046
$ws.Range("F33").Value = 56
$ws.Range("G33").Value = 0
$ws.Range("H33").Value = 49

# Row 118
$ws.Range("A118").Value = "Consejo Danes para los Refugiados"
$ws.Range("B118").Value = 65
$ws.Range("C118").Value = 14
$ws.Range("D118").Value = 2
$ws.Range("E118").Value = 57
$ws.Range("G118").Value = 3
$ws.Range("H118").Value = 6

# Row 119
$ws.Range("A119").Value = "Mayotte"
$ws.Range("B119").Value = 63
$ws.Range("C119").Value = 13
$ws.Range("E119").Value = 63

# Row 120
$ws.Range("A120").Value = "Ruanda"
$ws.Range("B120").Value = 60
$ws.Range("C120").Value = 6
$ws.Range("E120").Value = 60

# Row 121
$ws.Range("A121").Value = "Kirguistan"
$ws.Range("C121").Value = 0
$ws.Range("D121").Value = 0
$ws.Range("E121").Value = 58
$ws.Range("G121").Value = 0
$ws.Range("H121").Value = 0

# Row 153
$ws.Range("A153").Value = "San Martin (Parte Francesa)"

# Row 154
$ws.Range("A154").Value = "Dominica"

# Row 159
$ws.Range("A159").Value = "Mozambique"
$ws.Range("C159").Value = 1

# Row 160
$ws.Range("A160").Value = "Laos"
$ws.Range("C160").Value = 2

# Row 162
$ws.Range("A162").Value = "Surinam"
$ws.Range("C162").Value = 0

# Row 164
$ws.Range("A164").Value = "Haiti"
$ws.Range("C164").Value = 0

# Row 165
$ws.Range("A165").Value = "Islas Caimanes"
$ws.Range("C165").Value = 0

# Row 166
$ws.Range("A166").Value = "Guyana"
$ws.Range("C166").Value = 3

# Row 170
$ws.Range("A170").Value = "Granada"

# Row 171
$ws.Range("A171").Value = "Seychelles"

# Row 172
$ws.Range("A172").Value = "Gabon"
$ws.Range("C172").Value = 0

# Row 173
$ws.Range("A173").Value = "Zimbabue"
$ws.Range("C173").Value = 2

# Row 174
$ws.Range("A174").Value = "Eritrea"
$ws.Range("C174").Value = 0

# Row 175
$ws.Range("A175").Value = "Benin"

# Row 176
$ws.Range("A176").Value = "Santa Sede"
$ws.Range("C176").Value = 2

# Row 177
$ws.Range("A177").Value = "Mauritania"
$ws.Range("C177").Value = 2

# Row 181
$ws.Range("A181").Value = "San Bartolome"
$ws.Range("C181").Value = 0

# Row 182
$ws.Range("A182").Value = "Cabo Verde"
$ws.Range("C182").Value = 0

# Row 183
$ws.Range("A183").Value = "Nepal"
$ws.Range("C183").Value = 1
$ws.Range("D183").Value = 1
$ws.Range("H183").Value = 0

# Row 184
$ws.Range("A184").Value = "Sudan"
$ws.Range("C184").Value = 2
$ws.Range("D184").Value = 0
$ws.Range("H184").Value = 1

# Row 186
$ws.Range("A186").Value = "Islas Turcas y Caicos"
$ws.Range("C186").Value = 2

# Row 187
$ws.Range("A187").Value = "Angola"
$ws.Range("E187").Value = 4
$ws.Range("H187").Value = 0

# Row 188
$ws.Range("A188").Value = "Nicaragua"
$ws.Range("B188").Value = 4
$ws.Range("C188").Value = 0
$ws.Range("H188").Value = 1

# Row 189
$ws.Range("A189").Value = "Libia"
$ws.Range("C189").Value = 2

# Row 190
$ws.Range("A190").Value = "San Martin (Parte Holandesa)"

# Row 191
$ws.Range("A191").Value = "Republica de Africa Central"

# Row 192
$ws.Range("A192").Value = "Republica del Chad"

# Row 193
$ws.Range("A193").Value = "Butan"

# Row 194
$ws.Range("A194").Value = "Liberia"

# Row 195
$ws.Range("A195").Value = "Somalia"
$ws.Range("D195").Value = 0
$ws.Range("E195").Value = 3

# Row 196
$ws.Range("A196").Value = "Santa Lucia"
$ws.Range("D196").Value = 1
$ws.Range("H196").Value = 0

# Row 197
$ws.Range("A197").Value = "Gambia"
$ws.Range("B197").Value = 3
$ws.Range("H197").Value = 1

# Row 198
$ws.Range("A198").Value = "Anguila"

# Row 199
$ws.Range("A199").Value = "Guinea-Bisau"

# Row 200
$ws.Range("A200").Value = "San Cristobal y Nieves"

# Row 201
$ws.Range("A201").Value = "Islas Virgenes Britanicas"

# Row 202
$ws.Range("A202").Value = "Belice"
